$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: Grafica / Creati asset per il giardino / 1.5h / Sfondo e panchina
$ws.Cells.Item(16, 2).Value = "Grafica"
$ws.Cells.Item(16, 3).Value = "Creati asset per il giardino"
$ws.Cells.Item(16, 4).Value = 1.5/24
$ws.Cells.Item(16, 5).Value = "Sfondo e panchina"

# Row 17: 24/12/2024 / Programmazione / Ink / 1h / Settaggio storia spettro Uno
$ws.Cells.Item(17, 1).Value = 45650
$ws.Cells.Item(17, 2).Value = "Programmazione"
$ws.Cells.Item(17, 3).Value = "Ink"
$ws.Cells.Item(17, 4).Value = 1/24
$ws.Cells.Item(17, 5).Value = "Settaggio storia spettro Uno"

# Row 18: Scrittura / Spettro Uno / 2h / Tematica, intro, oggetti, domande
$ws.Cells.Item(18, 2).Value = "Scrittura"
$ws.Cells.Item(18, 3).Value = "Spettro Uno"
$ws.Cells.Item(18, 4).Value = 2/24
$ws.Cells.Item(18, 5).Value = "Tematica, intro, oggetti, domande"

[void]$ws.Range("E26").Select()
